$d = $word.ActiveDocument

# 1. Rubric (60 pts total) -> Rubric (100 pts total)
$d.Content.Find.Execute("Rubric (6", $false, $false, $false, $false, $false, $true, 1, $false, "Rubric (10", 2) | Out-Null

# 2. Title and author (5 pts) -> (10 pts)
$d.Content.Find.Execute("Title and author (5", $false, $false, $false, $false, $false, $true, 1, $false, "Title and author (10", 2) | Out-Null

# 3. Introduction (15 pts) -> (30 pts)
$d.Content.Find.Execute("Introduction (15", $false, $false, $false, $false, $false, $true, 1, $false, "Introduction (30", 2) | Out-Null

# 4. Description of problem - 5 pts -> - 10 pts
$dash = [char]8211
$d.Content.Find.Execute("Description of problem " + $dash + " 5 pts", $false, $false, $false, $false, $false, $true, 1, $false, "Description of problem " + $dash + " 10 pts", 2) | Out-Null

# 5. Scientific question - 5 pts -> - 10 pts
$d.Content.Find.Execute("Scientific question " + $dash + " 5 pts", $false, $false, $false, $false, $false, $true, 1, $false, "Scientific question " + $dash + " 10 pts", 2) | Out-Null

# 6. Background information with references - 5 pts -> - 10 pts
$d.Content.Find.Execute("Background information with references " + $dash + " 5 pts", $false, $false, $false, $false, $false, $true, 1, $false, "Background information with references " + $dash + " 10 pts", 2) | Out-Null

# 7. Methods (35 pts) -> (50 pts)
$d.Content.Find.Execute("Methods (35", $false, $false, $false, $false, $false, $true, 1, $false, "Methods (50", 2) | Out-Null

# 8. Background information on the field site (5 pts) -> (10 pts)
$d.Content.Find.Execute("Background information on the field site (5", $false, $false, $false, $false, $false, $true, 1, $false, "Background information on the field site (10", 2) | Out-Null

# 9. Map of sampling locations (5 pts) -> (10 pts)
$d.Content.Find.Execute("Map of sampling locations (5", $false, $false, $false, $false, $false, $true, 1, $false, "Map of sampling locations (10", 2) | Out-Null

# 10. Field collection (10 pts) -> (20 pts)
$d.Content.Find.Execute("Field collection (10", $false, $false, $false, $false, $false, $true, 1, $false, "Field collection (20", 2) | Out-Null

# 11. Description of sample locations (2.5 pts) -> (5 pts)
$d.Content.Find.Execute("Description of s" + "ample locations (2.5", $false, $false, $false, $false, $false, $true, 1, $false, "Description of s" + "ample locations (5", 2) | Out-Null

# 12. Sample collection, including equipment (2.5 pts) -> (5 pts)
$d.Content.Find.Execute("Sample collection, including equipment (2.5", $false, $false, $false, $false, $false, $true, 1, $false, "Sample collection, including equipment (5", 2) | Out-Null

# 13. Description of subsampling for different analyses (2.5 pts) -> (5 pts)
$d.Content.Find.Execute("Description of subsampling for different analyses (2.5", $false, $false, $false, $false, $false, $true, 1, $false, "Description of subsampling for different analyses (5", 2) | Out-Null

# 14. Field measurements, including equipment (2.5 pts) -> (5 pts)
$d.Content.Find.Execute("Field measurements, including equipment (2.5", $false, $false, $false, $false, $false, $true, 1, $false, "Field measurements, including equipment (5", 2) | Out-Null

# 15. Laboratory analysis (10 pts) -> (5 pts)
$d.Content.Find.Execute("Laboratory analysis (10", $false, $false, $false, $false, $false, $true, 1, $false, "Laboratory analysis (5", 2) | Out-Null

# 16. Table of samples with field measurements (5 pts) -> (10 pts)
$d.Content.Find.Execute("Table of samples with field measurements (5", $false, $false, $false, $false, $false, $true, 1, $false, "Table of samples with field measurements (10", 2) | Out-Null

for ($i = 1; $i -le $d.Paragraphs.Count; $i++) {
    Write-Output ("$i : " + $d.Paragraphs.Item($i).Range.Text)
}
